$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1413
$ws.Range("I2").Value = 3798
$ws.Range("J2").Value = 15507
$ws.Range("K2").Value = 71
$ws.Range("L2").Value = 4241
$ws.Range("M2").Value = 245
$ws.Range("N2").Value = 2651
$ws.Range("O2").Value = 16
$ws.Range("Q2").Value = 36
$ws.Range("R2").Value = 197
$ws.Range("S2").Value = 1660
$ws.Range("T2").Value = 2684
$ws.Range("U2").Value = 237
$ws.Range("V2").Value = 23701
$ws.Range("X2").Value = 24067
$ws.Range("Y2").Value = 25
$ws.Range("Z2").Value = 353
$ws.Range("AA2").Value = 150
